# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column as plain text so values such as "1.003" or
# "23.876.75" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Rows 37 and 38 swapped rank (FraxShare <-> InternetComputer(DFINITY))
# along with their own updated Price / Volume(1h) values.
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "6.648"
$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "10.90"
$ws.Range("E38").Value = "  +5.34%  "

# Update Price (D) and Volume(1h) (E) columns for all other rows whose
# values changed in this data refresh.
$ws.Range("D2").Value = "23.879.35"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.654.11"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "309.38"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "0.3893"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").Value = "0.3839"
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "50.93"
$ws.Range("D10").Value = "1.351"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "0.08461"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "23.85"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "7.162"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "7.940"
$ws.Range("E15").Value = "  +7.55%  "
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("D17").Value = "1.656.74"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "94.40"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "0.06993"
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "19.81"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "6.924"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "13.62"
$ws.Range("E23").Value = "  +2.09%  "
$ws.Range("D24").Value = "23.876.75"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "2.479"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("D26").Value = "3.042"
$ws.Range("E26").Value = "  +10.10%  "
$ws.Range("D27").Value = "22.09"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "152.99"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").Value = "5.359"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("D30").Value = "139.42"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "7.840"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").Value = "2.494"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "1.838.64"
$ws.Range("E33").Value = "  +3.17%  "
$ws.Range("D34").Value = "1.039"
$ws.Range("E34").Value = "  +8.41%  "
$ws.Range("D35").Value = "0.08068"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "0.02964"
$ws.Range("E36").Value = "  +4.27%  "
$ws.Range("D39").Value = "0.2685"
$ws.Range("E39").Value = "  +1.74%  "
$ws.Range("D40").Value = "0.09125"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "13.49"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").Value = "0.7516"
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("D43").Value = "1.416"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").Value = "16.26"
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("D45").Value = "0.6958"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("D46").Value = "2.461"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("D49").Value = "0.08246"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").Value = "134.41"
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("D51").Value = "1.232"
$ws.Range("E51").Value = "  -0.50%  "
